# Add the new customer row (row 64): phone 51616176, no birthday on file,
# 0 total_points — mirrors the existing rows' layout where the phone
# number is stored as text and a missing birthday is an empty text cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Leading "'" forces these to be entered as text (matching the sheet's
# existing phone-number / blank-birthday cells) instead of Excel's default
# "looks like a number" auto-typing.
$ws.Cells.Item($row, 1).Value = "'51616176"
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 3).Value = 0

# The quote-prefix entry above picks up a "quoted text" cell style; put the
# cells back on the sheet's default (unstyled) look, same as the rest of
# the data rows.
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
